# DataEngine.xlsx update
# Adds new "verifyElement" and "compareLinkText" action-keyword test steps
# (Verify Username / Password / Login button webelements, and a new
# "Verify_Account" test case block) to the "Test Steps" sheet, and a
# matching row on the "Test Cases" sheet. Also fixes up column widths,
# the active sheet/selection, and the hyperlink cell references that
# moved as a result of the new rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Test Steps"
$ws2 = $wb.Worksheets.Item(2)   # "Test Cases"

# ---------------------------------------------------------------------
# 1. Make room for the new "Test Steps" rows by inserting blank rows at
#    the right spots (this shifts everything below down and keeps the
#    existing row formatting on rows that move).
# ---------------------------------------------------------------------
$ws1.Rows.Item(4).Insert()   # new TS_003 "Verify Username" (before old Enter-Email row)
$ws1.Rows.Item(6).Insert()   # new TS_005 "Verify Password" (before old Enter-password row)
$ws1.Rows.Item(8).Insert()   # new TS_007 "Verify login button" (before old Click-Signin row)
$ws1.Rows.Item(11).Insert()  # new Verify_Account block, row 1
$ws1.Rows.Item(11).Insert()  # new Verify_Account block, row 2

# Give the newly inserted (currently blank) rows the same borders/style
# used by the rest of the data rows (copy format from row 2).
$ws1.Range("A2:G2").Copy($ws1.Range("A4:G4"))
$ws1.Range("A2:G2").Copy($ws1.Range("A6:G6"))
$ws1.Range("A2:G2").Copy($ws1.Range("A8:G8"))
$ws1.Range("A2:G2").Copy($ws1.Range("A11:G11"))
$ws1.Range("A2:G2").Copy($ws1.Range("A12:G12"))

# ---------------------------------------------------------------------
# 2. Re-write every TS_ID cell (column B) so the sequence is contiguous
#    TS_001 .. TS_015 again, and fill in the new rows' content.
# ---------------------------------------------------------------------
$ws1.Range("B2").Value2  = "TS_001"
$ws1.Range("B3").Value2  = "TS_002"

$ws1.Range("B4").Value2  = "TS_003"
$ws1.Range("C4").Value2  = "Verify the webelement Username on page"
$ws1.Range("D4").Value2  = "txt_Username"
$ws1.Range("E4").Value2  = "verifyElement"
$ws1.Range("F4").Value2  = ""

$ws1.Range("B5").Value2  = "TS_004"

$ws1.Range("B6").Value2  = "TS_005"
$ws1.Range("C6").Value2  = "Verify the webelement Passwrod on page"
$ws1.Range("D6").Value2  = "txt_Password"
$ws1.Range("E6").Value2  = "verifyElement"
$ws1.Range("F6").Value2  = ""

$ws1.Range("B7").Value2  = "TS_006"

$ws1.Range("B8").Value2  = "TS_007"
$ws1.Range("C8").Value2  = "Verify the webelement login button"
$ws1.Range("D8").Value2  = "btn_Login"
$ws1.Range("E8").Value2  = "verifyElement"
$ws1.Range("F8").Value2  = ""

$ws1.Range("B9").Value2  = "TS_008"
$ws1.Range("B10").Value2 = "TS_009"

# New "Verify_Account" test-case block (rows 11-12)
$ws1.Range("A11").Value2 = "Verify_Account"
$ws1.Range("B11").Value2 = "TS_010"
$ws1.Range("C11").Value2 = "Compare link text of the customer account who has logged-in"
$ws1.Range("D11").Value2 = "lnk_custaccount"
$ws1.Range("E11").Value2 = "compareLinkText"
$ws1.Range("F11").Value2 = "ravi reddy"

$ws1.Range("A12").Value2 = "Verify_Account"
$ws1.Range("B12").Value2 = "TS_011"
$ws1.Range("C12").Value2 = "Wait for sometime"
$ws1.Range("D12").Value2 = ""
$ws1.Range("E12").Value2 = "waitFor"
$ws1.Range("F12").Value2 = ""

# Existing Menu_Navigation / Logout_01 rows, renumbered
$ws1.Range("B13").Value2 = "TS_012"
$ws1.Range("B14").Value2 = "TS_013"
$ws1.Range("B15").Value2 = "TS_014"
$ws1.Range("B16").Value2 = "TS_015"

# ---------------------------------------------------------------------
# 3. Hyperlinks on column F moved from F4/F5 to F5/F7 - drop the old
#    ones and add fresh hyperlinks at the new locations.
# ---------------------------------------------------------------------
$ws1.Hyperlinks.Delete()
$ws1.Range("F5").Value2 = "ravikaanthe@rediffmail.com"
$ws1.Range("F7").Value2 = "test@123"
$ws1.Hyperlinks.Add($ws1.Range("F5"), "mailto:ravikaanthe@rediffmail.com") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("F7"), "mailto:test@123") | Out-Null

# ---------------------------------------------------------------------
# 4. Column width tweaks (columns C and D got a bit wider to fit the
#    new, longer descriptions/page-object names).
# ---------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 51.9
$ws1.Columns.Item(4).ColumnWidth = 14.3

# ---------------------------------------------------------------------
# 5. "Test Cases" sheet - add the matching "Verify_Account" row, which
#    is inserted right after "Login_01" (so Menu_Navigation/Logout_01
#    shift down by one row).
# ---------------------------------------------------------------------
$ws2.Rows.Item(3).Insert()
$ws2.Range("A2:D2").Copy($ws2.Range("A3:D3"))

$ws2.Range("A3").Value2 = "Verify_Account"
$ws2.Range("B3").Value2 = "Verifying the account of user who logged-in"
$ws2.Range("C3").Value2 = "Yes"
$ws2.Range("D3").Value2 = "PASS"

# ---------------------------------------------------------------------
# 6. Sheet/selection state: "Test Steps" becomes the active tab (instead
#    of "Test Cases"), with new cell selections on each sheet.
# ---------------------------------------------------------------------
$ws2.Range("B3").Select()
$ws1.Select()
$ws1.Range("E2").Select()
